$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new parameter row (A52/B52): "latticehold" with unit "ms" ---
$ws.Range("A52").Value = "latticehold"
$ws.Range("B52").Value = "ms"

# --- Adjust column widths (stored OOXML width: col A 16->17, col B 19.3984375->~20.125) ---
$ws.Columns.Item(1).ColumnWidth = 16.22
$ws.Columns.Item(2).ColumnWidth = 19.4

# --- Request a full recalculation on next load ---
$wb.ForceFullCalculation = $true
